$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.129.40"
$ws.Range("E2").Value = "  +0.00%  "

$ws.Range("D3").Value = "2.047.61"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").Value = "248.59"
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("D6").Value = "0.666"
$ws.Range("E6").Value = "  -1.36%  "

$ws.Range("D7").Value = "59.51"
$ws.Range("E7").Value = "  +2.64%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "0.384"
$ws.Range("E9").Value = "  +1.18%  "

$ws.Range("D10").Value = "0.0787"
$ws.Range("E10").Value = "  -2.87%  "

$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("D12").Value = "15.85"
$ws.Range("E12").Value = "  +4.27%  "

$ws.Range("D13").Value = "2.338.92"
$ws.Range("E13").Value = "  -0.87%  "

$ws.Range("D14").Value = "0.835"
$ws.Range("E14").Value = "  +2.20%  "

$ws.Range("D15").Value = "5.75"
$ws.Range("E15").Value = "  +7.62%  "

$ws.Range("D16").Value = "2.053.50"
$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("D17").Value = "18.28"
$ws.Range("E17").Value = "  +26.38%  "

$ws.Range("D18").Value = "37.102.02"
$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("D19").Value = "74.88"
$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("D20").Value = "0.0₃0896"
$ws.Range("E20").Value = "  -2.70%  "

$ws.Range("D21").Value = "5.35"
$ws.Range("E21").Value = "  -0.47%  "

$ws.Range("D22").Value = "237.21"
$ws.Range("E22").Value = "  -0.64%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").Value = "2.44"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "2.17"
$ws.Range("E25").Value = "  +7.92%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "169.04"
$ws.Range("E26").Value = "  -1.66%  "

$ws.Range("D27").Value = "9.39"
$ws.Range("E27").Value = "  +2.39%  "

$ws.Range("D28").Value = "20.05"
$ws.Range("E28").Value = "  -0.89%  "

$ws.Range("D29").Value = "0.125"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").Value = "1.14"
$ws.Range("E30").Value = "  +5.04%  "

$ws.Range("D31").Value = "4.78"
$ws.Range("E31").Value = "  +3.29%  "

$ws.Range("D32").Value = "0.0626"

$ws.Range("D33").Value = "4.54"
$ws.Range("E33").Value = "  +2.58%  "

$ws.Range("D34").Value = "0.0893"
$ws.Range("E34").Value = "  +0.65%  "

$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("D36").Value = "2.21"
$ws.Range("E36").Value = "  -2.72%  "

$ws.Range("D37").Value = "1.75"
$ws.Range("E37").Value = "  -1.42%  "

$ws.Range("D38").Value = "1.34"
$ws.Range("E38").Value = "  -0.57%  "

$ws.Range("D39").Value = "0.106"
$ws.Range("E39").Value = "  -3.78%  "

$ws.Range("D40").Value = "3.17"
$ws.Range("E40").Value = "  +13.76%  "

$ws.Range("D41").Value = "5.18"
$ws.Range("E41").Value = "  +15.91%  "

$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.50"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.99%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0221"
$ws.Range("E43").Value = "  -1.73%  "

$ws.Range("D44").Value = "1.13"
$ws.Range("E44").Value = "  -1.54%  "

$ws.Range("D45").Value = "96.04"
$ws.Range("E45").Value = "  -1.06%  "

$ws.Range("D46").Value = "2.47"
$ws.Range("E46").Value = "  -0.95%  "

$ws.Range("D47").Value = "2.91"
$ws.Range("E47").Value = "  -0.35%  "

$ws.Range("D48").Value = "1.281.91"
$ws.Range("E48").Value = "  -1.88%  "

$ws.Range("D49").Value = "6.77"
$ws.Range("E49").Value = "  -1.84%  "

$ws.Range("D50").Value = "2.220.92"
$ws.Range("E50").Value = "  -1.04%  "

$ws.Range("D51").Value = "44.02"
$ws.Range("E51").Value = "  +0.82%  "
